$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set up I1/IF headers with the same style as the other header cells (copy format from H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows: row index, I0 value, IF value
$data = @(
    @(2, 7, 8),
    @(3, 4, 4),
    @(4, 6, 7),
    @(5, 8, 8),
    @(6, 8, 8),
    @(7, 7, 7),
    @(8, 8, 8),
    @(9, 6, 6),
    @(10, 9, 9),
    @(11, 7, 7),
    @(12, 8, 8),
    @(13, 5, 6),
    @(14, 7, 7),
    @(15, 6, 6),
    @(16, 6, 6),
    @(17, 8, 8),
    @(18, 5, 6),
    @(19, 7, 7),
    @(20, 7, 7),
    @(21, 6, 6),
    @(22, 8, 8),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 8, 8),
    @(26, 11, 11),
    @(27, 9, 9),
    @(28, 6, 6),
    @(29, 7, 7),
    @(30, 5, 5),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 9, 9),
    @(34, 6, 7),
    @(35, 8, 8),
    @(36, 7, 7),
    @(37, 7, 7),
    @(38, 6, 6),
    @(39, 8, 8),
    @(40, 8, 8),
    @(41, 8, 8),
    @(42, 6, 6),
    @(43, 7, 7),
    @(44, 7, 7),
    @(45, 6, 7),
    @(46, 9, 9),
    @(47, 6, 7),
    @(48, 7, 7),
    @(49, 7, 7),
    @(50, 8, 8),
    @(51, 8, 8),
    @(52, 7, 8),
    @(53, 9, 9),
    @(54, 7, 8),
    @(55, 8, 8),
    @(56, 8, 8),
    @(57, 6, 7),
    @(58, 9, 9),
    @(59, 7, 8),
    @(60, 8, 8),
    @(61, 8, 8),
    @(62, 6, 7),
    @(63, 8, 9),
    @(64, 9, 9),
    @(65, 8, 8),
    @(66, 9, 9),
    @(67, 7, 7),
    @(68, 6, 7),
    @(69, 6, 7),
    @(70, 7, 8),
    @(71, 7, 7),
    @(72, 7, 7),
    @(73, 7, 7),
    @(74, 6, 6),
    @(75, 7, 7),
    @(76, 10, 10),
    @(77, 8, 8),
    @(78, 6, 7),
    @(79, 6, 6),
    @(80, 6, 6),
    @(81, 5, 5),
    @(82, 6, 6),
    @(83, 5, 5),
    @(84, 2, 2),
    @(85, 4, 4)
)

foreach ($item in $data) {
    $r = $item[0]
    $iVal = $item[1]
    $jVal = $item[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
